$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New case entries for 21TRD09437 / Bunner now that dismissed-charge rows are
# tracked with their own "Dismissed" row - the guilty buttons need blocking
# once a case is fully resolved, so the three underlying counts are written
# out individually (DUS, Reckless Operation, and the dismissed minor
# misdemeanor charge).

# Row 85: DUS / No Contest / Guilty
$ws.Range("A85").Value = "21TRD09437"
$ws.Range("B85").Value = "Bunner"
$ws.Range("C85").Value = "DUS"
$ws.Range("D85").Value = "'4510.11"
$ws.Range("D85").ClearFormats()
$ws.Range("E85").Value = "M1"
$ws.Range("F85").Value = "No Contest"
$ws.Range("G85").Value = "Guilty"
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = "'0"
$ws.Range("I85").ClearFormats()

# Row 86: RECKLESS OPERATION 1ST IN 1 YR / No Contest / Guilty - Allied Offense
$ws.Range("A86").Value = "21TRD09437"
$ws.Range("B86").Value = "Bunner"
$ws.Range("C86").Value = "RECKLESS OPERATION 1ST IN 1 YR"
$ws.Range("D86").Value = "'4511.20"
$ws.Range("D86").ClearFormats()
$ws.Range("E86").Value = "MM"
$ws.Range("F86").Value = "No Contest"
$ws.Range("G86").Value = "Guilty - Allied Offense"
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = "'0"
$ws.Range("I86").ClearFormats()

# Row 87: blank offense/statute (empty text cells), Minor Misdemeanor, Dismissed / Guilty
$ws.Range("A87").Value = "21TRD09437"
$ws.Range("B87").Value = "Bunner"
$ws.Range("C87").Value = "'"
$ws.Range("C87").ClearFormats()
$ws.Range("D87").Value = "'"
$ws.Range("D87").ClearFormats()
$ws.Range("E87").Value = "Minor Misdemeanor"
$ws.Range("F87").Value = "Dismissed"
$ws.Range("G87").Value = "Guilty"
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = "'0"
$ws.Range("I87").ClearFormats()
